$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
}

Set-TextValue "D2" "51.894.07"
Set-TextValue "E2" "  +1.65%  "
Set-TextValue "D3" "2.813.15"
Set-TextValue "E3" "  +1.97%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "351.47"
Set-TextValue "E5" "  -0.36%  "
Set-TextValue "D6" "113.04"
Set-TextValue "E6" "  +4.90%  "
Set-TextValue "D7" "0.558"
Set-TextValue "E7" "  +1.89%  "
Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  -0.05%  "
Set-TextValue "D9" "0.622"
Set-TextValue "E9" "  +6.50%  "
Set-TextValue "D10" "40.46"
Set-TextValue "E10" "  +2.70%  "
Set-TextValue "D12" "0.0841"
Set-TextValue "E12" "  +0.84%  "
Set-TextValue "D13" "19.89"
Set-TextValue "E13" "  +0.47%  "
Set-TextValue "D14" "7.79"
Set-TextValue "E14" "  +3.86%  "
Set-TextValue "D15" "3.248.19"
Set-TextValue "E15" "  +1.87%  "
Set-TextValue "D16" "0.970"
Set-TextValue "E16" "  +4.29%  "
Set-TextValue "D17" "2.814.35"
Set-TextValue "E17" "  +1.85%  "
Set-TextValue "D18" "51.897.04"
Set-TextValue "E18" "  +1.79%  "
Set-TextValue "D19" "3.38"
Set-TextValue "E19" "  +10.22%  "
Set-TextValue "D20" "7.66"
Set-TextValue "E20" "  -0.07%  "
Set-TextValue "D21" "13.62"
Set-TextValue "E21" "  +4.46%  "
Set-TextValue "D22" "0.0₃0976"
Set-TextValue "E22" "  +1.72%  "
Set-TextValue "D23" "70.51"
Set-TextValue "E23" "  +1.34%  "
Set-TextValue "D24" "269.64"
Set-TextValue "E24" "  +1.79%  "
Set-TextValue "E25" "  +1.87%  "
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  -0.05%  "
Set-TextValue "D27" "26.22"
Set-TextValue "E27" "  +1.26%  "
Set-TextValue "D28" "0.163"
Set-TextValue "E28" "  +1.08%  "
Set-TextValue "D29" "38.98"
Set-TextValue "E29" "  +13.42%  "
Set-TextValue "D30" "10.48"
Set-TextValue "E31" "  +1.01%  "
Set-TextValue "D32" "52.78"
Set-TextValue "E32" "  +2.16%  "
Set-TextValue "D33" "6.16"
Set-TextValue "E33" "  +1.77%  "
Set-TextValue "D34" "0.0904"
Set-TextValue "E34" "  +9.05%  "
Set-TextValue "E35" "  +2.70%  "
Set-TextValue "D36" "5.67"
Set-TextValue "E36" "  +4.53%  "
Set-TextValue "D37" "0.999"
Set-TextValue "E37" "  -0.01%  "
Set-TextValue "D38" "18.97"
Set-TextValue "E38" "  +3.69%  "
Set-TextValue "D39" "3.20"
Set-TextValue "E39" "  +2.06%  "
Set-TextValue "D40" "2.02"
Set-TextValue "E40" "  +3.22%  "
Set-TextValue "E41" "  +2.33%  "
Set-TextValue "D42" "2.53"
Set-TextValue "E42" "  +1.11%  "
Set-TextValue "D43" "122.00"
Set-TextValue "E43" "  +1.37%  "
Set-TextValue "E44" "  +1.98%  "
Set-TextValue "D45" "22.04"
Set-TextValue "E45" "  -1.11%  "
Set-TextValue "D46" "3.55"
Set-TextValue "E46" "  +9.86%  "
Set-TextValue "D47" "2.47"
Set-TextValue "E47" "  +9.40%  "
Set-TextValue "D48" "2.127.88"
Set-TextValue "E48" "  +2.21%  "
Set-TextValue "D49" "0.997"
Set-TextValue "E49" "  +8.75%  "
Set-TextValue "D50" "5.51"
Set-TextValue "E50" "  +0.71%  "
Set-TextValue "D51" "0.223"
Set-TextValue "E51" "  +17.42%  "
